$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 4.4
$ws.Range("L2").Value = 1.31
$ws.Range("N2").Value = 4
$ws.Range("P2").Value = 2.04
$ws.Range("R2").Value = 1.4
$ws.Range("U2").Value = 2
$ws.Range("AF2").Value = 11.5
$ws.Range("AG2").Value = 10
$ws.Range("AI2").Value = 1000
$ws.Range("AN2").Value = 9.800000000000001
$ws.Range("F3").Value = 6.4
$ws.Range("G3").Value = 8.4
$ws.Range("H3").Value = 1.46
$ws.Range("K3").Value = 5.3
$ws.Range("Q3").Value = 1.68
$ws.Range("F5").Value = 1.4
$ws.Range("S5").Value = 3.3
$ws.Range("T5").Value = 2.24
$ws.Range("F6").Value = 2.12
$ws.Range("G6").Value = 2.82
$ws.Range("I6").Value = 4.2
$ws.Range("J6").Value = 2.96
$ws.Range("K6").Value = 5.4
$ws.Range("F7").Value = 2.66
$ws.Range("G7").Value = 2.7
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3.05
$ws.Range("I8").Value = 1.62
$ws.Range("G9").Value = 4.7
$ws.Range("I9").Value = 2.14
$ws.Range("Q9").Value = 2.3
$ws.Range("G10").Value = 3.85
$ws.Range("H10").Value = 2.38
$ws.Range("Q10").Value = 2.46
$ws.Range("G11").Value = 2.32
$ws.Range("H11").Value = 3.35
$ws.Range("I11").Value = 4.9
$ws.Range("J11").Value = 3.4
$ws.Range("K11").Value = 6.8
$ws.Range("F12").Value = 2.38
$ws.Range("G12").Value = 2.4
$ws.Range("H12").Value = 3.35
$ws.Range("J13").Value = 3.8
$ws.Range("S13").Value = 3.75
$ws.Range("X13").Value = 13.5
$ws.Range("AB13").Value = 8
$ws.Range("F14").Value = 3.55
$ws.Range("P14").Value = 1.78
$ws.Range("Q14").Value = 2.08
$ws.Range("P15").Value = 2.14
$ws.Range("AK15").Value = 46
$ws.Range("H16").Value = 21
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 9.199999999999999
$ws.Range("K16").Value = 11
$ws.Range("Q16").Value = 1.39
